$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the FFR (C) column regression coefficients for the lag rows ---
# "-0.0" and "-0.005" look like plain numbers to Excel's auto-detection, so
# force them to be typed in as text (matching the original file, where every
# coefficient cell is stored as text) and then drop the resulting cell
# format so no stray "Text" number-format style is left behind.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "-0.0"
$ws.Range("C2").ClearFormats()

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "-0.005"
$ws.Range("C3").ClearFormats()

$ws.Range("C4").Value = "0.069**"

# --- Update the C column regression coefficients for the lag rows ---
$ws.Range("D2").Value = "-1.814**"
$ws.Range("D3").Value = "-2.304***"
$ws.Range("D4").Value = "-0.663***"

# --- Remove the "Constant" and "r2_adj" rows entirely ---
$ws.Range("A5:D6").Delete()
